# Commit: feat: add 2022-Q4 data
#
# A new quarterly snapshot sheet "2022-Q4" is inserted right after the
# "总计" (summary) sheet, pushing all the existing quarter sheets one
# position to the right (their own names/content are untouched). The
# "总计" sheet gets a new first data row summarising the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new worksheet right after "总计" (i.e. before "2022-Q3")
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $zongji)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2. Populate "2022-Q4" with the fund holdings table (same layout used
#    by every other quarterly sheet in this workbook).
# ---------------------------------------------------------------------
$q4.Range("B1:H1").NumberFormat = "@"
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Row 2: 005051 - 上投摩根标普港股通低波红利指数A
$q4.Range("A2").Value = 0
$q4.Range("B2:G2").NumberFormat = "@"
$q4.Range("B2").Value = "005051"
$q4.Range("C2").Value = "上投摩根标普港股通低波红利指数A"
$q4.Range("D2").Value = "1.64"
$q4.Range("E2").Value = "93.98"
$q4.Range("F2").Value = "2.85"
$q4.Range("G2").Value = "0.0467"
$q4.Range("H2").Value = 3

# Row 3: 005052 - 上投摩根标普港股通低波红利指数C
$q4.Range("A3").Value = 1
$q4.Range("B3:G3").NumberFormat = "@"
$q4.Range("B3").Value = "005052"
$q4.Range("C3").Value = "上投摩根标普港股通低波红利指数C"
$q4.Range("D3").Value = "1.52"
$q4.Range("E3").Value = "93.98"
$q4.Range("F3").Value = "2.85"
$q4.Range("G3").Value = "0.0433"
$q4.Range("H3").Value = 3

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: insert a new row 2 for 2022-Q4 and
#    shift the other quarters' totals down (their B/C/D values stay the
#    same, only their row position and running index in column A move).
# ---------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q4"
$zongji.Range("C2").Value = 2
$zongji.Range("D2").Value = 0.09

# Renumber the running index in column A for the rows that shifted down
# (rows 3..9 now hold what used to be rows 2..8, index 0..6 -> 1..7).
for ($r = 3; $r -le 9; $r++) {
    $zongji.Cells.Item($r, 1).Value = ($r - 2)
}

Write-Host "2022-Q4 sheet added and 总计 updated"
